$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = "<rpc-reply message-id=""urn:uuid:f693dd02-ccd4-4ec2-af38-a8eff117569b"">`n  <data/>`n</rpc-reply>`n"

$ws.Range("G2").Value = "  <edit-config>`n    <target>`n      <candidate/>`n    </target>`n    <config>`n      <network-instances xmlns=""http://openconfig.net/yang/network-instance"">`n        <network-instance>`n          <name>Prueba_LxVPN</name>`n          <config>`n              <name>Prueba_LxVPN</name>`n              <type xmlns:oc-ni-types=""http://openconfig.net/yang/network-instance-types"">oc-ni-types:L3VRF</type>`n              <enabled-address-families xmlns:oc-types=""http://openconfig.net/yang/openconfig-types"">oc-types:IPV4</enabled-address-families>`n          </config>`n          <encapsulation>`n            <config>`n              <label-allocation-mode xmlns:oc-ni-types=""http://openconfig.net/yang/network-instance-types"">oc-ni-types:PER_PREFIX</label-allocation-mode>`n            </config>`n          </encapsulation>`n        </network-instance>`n      </network-instances>`n    </config>`n  </edit-config>"

$ws.Range("H2").Value = "- Response of edit-config: <rpc-reply message-id=""urn:uuid:7e56925c-4215-49db-babf-47161cc92757"">`n  <ok/>`n</rpc-reply>`n - Response of commit: <rpc-reply xmlns:nc-ext=""urn:huawei:yang:huawei-ietf-netconf-ext"" message-id=""urn:uuid:1c6a1bbb-ebca-4c03-8c20-57e7a4a42ae4"" nc-ext:flow-id=""249"">`n  <ok/>`n</rpc-reply>`n"

$ws.Range("I2").Value = "<rpc-reply message-id=""urn:uuid:9ac1e210-6ced-48ec-9859-ca9b9524ce01"">`n  <data>`n    <network-instances>`n      <network-instance>`n        <name>Prueba_LxVPN</name>`n        <config>`n          <name>Prueba_LxVPN</name>`n          <type>oc-ni-types:L3VRF</type>`n          <enabled-address-families>oc-types:IPV4</enabled-address-families>`n        </config>`n        <protocols>`n          <protocol>`n            <identifier>oc-pol-types:STATIC</identifier>`n            <name>default</name>`n            <config>`n              <identifier>oc-pol-types:STATIC</identifier>`n              <name>default</name>`n            </config>`n          </protocol>`n          <protocol>`n            <identifier>oc-pol-types:DIRECTLY_CONNECTED</identifier>`n            <name>default</name>`n            <config>`n              <identifier>oc-pol-types:DIRECTLY_CONNECTED</identifier>`n              <name>default</name>`n            </config>`n          </protocol>`n        </protocols>`n      </network-instance>`n    </network-instances>`n  </data>`n</rpc-reply>`n"
